# Updated cryptos list values (price/volume refresh), per the commit diff.
# Source data cells are plain text (t="inlineStr"/shared string); numeric-looking
# Price values are written with a leading apostrophe so Excel stores them as text
# (matching the original "99.99"/"x.xxx.xx" style strings) instead of converting them
# to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.133.61"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3
$ws.Range("D3").Value = "3.381.78"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'574.34"
$ws.Range("E5").Value = "  -0.94%  "

# Row 6
$ws.Range("D6").Value = "'137.48"
$ws.Range("E6").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").Value = "3.380.56"
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("E9").Value = "  -1.25%  "

# Row 10
$ws.Range("D10").Value = "'7.64"
$ws.Range("E10").Value = "  +1.85%  "

# Row 11
$ws.Range("E11").Value = "  -3.21%  "

# Row 12
$ws.Range("D12").Value = "'0.381"
$ws.Range("E12").Value = "  -2.71%  "

# Row 13
$ws.Range("D13").Value = "3.961.61"
$ws.Range("E13").Value = "  -0.23%  "

# Row 14
$ws.Range("E14").Value = "  +0.71%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'25.78"
$ws.Range("E15").Value = "  +1.34%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("E16").Value = "  -3.18%  "

# Row 17
$ws.Range("D17").Value = "3.385.84"
$ws.Range("E17").Value = "  -0.25%  "

# Row 18
$ws.Range("D18").Value = "61.264.63"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19
$ws.Range("D19").Value = "'13.83"
$ws.Range("E19").Value = "  -2.27%  "

# Row 20
$ws.Range("E20").Value = "  -1.23%  "

# Row 21
$ws.Range("D21").Value = "'9.34"
$ws.Range("E21").Value = "  -1.65%  "

# Row 22
$ws.Range("D22").Value = "'376.90"
$ws.Range("E22").Value = "  -1.53%  "

# Row 23
$ws.Range("B23").Value = "WrappedeETH"
$ws.Range("C23").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D23").Value = "3.518.50"
$ws.Range("E23").Value = "  -0.39%  "

# Row 24
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.551"
$ws.Range("E24").Value = "  -2.48%  "

# Row 25
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("E26").Value = "  -2.02%  "

# Row 27
$ws.Range("D27").Value = "'71.00"
$ws.Range("E27").Value = "  -0.50%  "

# Row 28
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'1.67"
$ws.Range("E28").Value = "  -3.18%  "

# Row 29
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.180"
$ws.Range("E29").Value = "  +11.84%  "

# Row 30
$ws.Range("E30").Value = "  -0.01%  "

# Row 31
$ws.Range("D31").Value = "'7.42"
$ws.Range("E31").Value = "  -3.11%  "

# Row 32
$ws.Range("D32").Value = "'8.06"
$ws.Range("E32").Value = "  -2.13%  "

# Row 33
$ws.Range("E33").Value = "  -1.71%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").Value = "'23.41"
$ws.Range("E35").Value = "  -0.38%  "

# Row 36
$ws.Range("E36").Value = "  -4.68%  "

# Row 37
$ws.Range("E37").Value = "  -2.09%  "

# Row 38
$ws.Range("E38").Value = "  -0.83%  "

# Row 39
$ws.Range("D39").Value = "'164.66"
$ws.Range("E39").Value = "  -0.53%  "

# Row 40
$ws.Range("D40").Value = "'0.0759"
$ws.Range("E40").Value = "  -3.75%  "

# Row 41
$ws.Range("D41").Value = "'25.58"
$ws.Range("E41").Value = "  +1.40%  "

# Row 42
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("E43").Value = "  -1.37%  "

# Row 44
$ws.Range("E44").Value = "  -2.81%  "

# Row 45
$ws.Range("D45").Value = "'4.34"

# Row 46
$ws.Range("E46").Value = "  -4.29%  "

# Row 47
$ws.Range("D47").Value = "2.546.44"
$ws.Range("E47").Value = "  +8.33%  "

# Row 48
$ws.Range("D48").Value = "'6.78"
$ws.Range("E48").Value = "  -1.54%  "

# Row 49
$ws.Range("D49").Value = "'22.88"
$ws.Range("E49").Value = "  -0.56%  "

# Row 50
$ws.Range("D50").Value = "'2.43"
$ws.Range("E50").Value = "  +3.78%  "

# Row 51
$ws.Range("D51").Value = "'0.0258"
$ws.Range("E51").Value = "  -1.83%  "

Write-Host "Updated cryptos list"
